$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 997.26
$ws.Range("I15").Value = 997.26
$ws.Range("K15").Value = 2991.78
$ws.Range("M15").Value = -2822.78
$ws.Range("H17").Value = 993.4643
$ws.Range("J17").Value = 1075.7084
$ws.Range("L17").Value = 3227.1252
$ws.Range("N17").Value = -3563.1252
$ws.Range("H53").Value = 400.375
$ws.Range("I53").Value = 167.8421
$ws.Range("J53").Value = 1284
$ws.Range("K53").Value = 167.8421
$ws.Range("L53").Value = 1284
$ws.Range("M53").Value = 469.1579
$ws.Range("N53").Value = -2558
$ws.Range("H70").Value = 1179.3871
$ws.Range("I70").Value = 1353.0869
$ws.Range("J70").Value = 680
$ws.Range("K70").Value = 4059.2607
$ws.Range("L70").Value = 2040
$ws.Range("M70").Value = -3789.2607
$ws.Range("N70").Value = -2580
$ws.Range("H73").Value = 1179.3871
$ws.Range("I73").Value = 1353.0869
$ws.Range("J73").Value = 680
$ws.Range("K73").Value = 4059.2607
$ws.Range("L73").Value = 2040
$ws.Range("M73").Value = -3123.2607
$ws.Range("N73").Value = -3912
$ws.Range("H120").Value = 38500
$ws.Range("J120").Value = 38500
$ws.Range("L120").Value = 38500
$ws.Range("N120").Value = -48176
$ws.Range("H129").Value = 1662.1708
$ws.Range("I129").Value = 564
$ws.Range("J129").Value = 2172.0356
$ws.Range("K129").Value = 1692
$ws.Range("L129").Value = 6516.1068
$ws.Range("M129").Value = 3308
$ws.Range("N129").Value = -16516.1068
$ws.Range("H132").Value = 5315.4756
$ws.Range("I132").Value = 4275.2173
$ws.Range("J132").Value = 8505.6
$ws.Range("K132").Value = 12825.6519
$ws.Range("L132").Value = 25516.8
$ws.Range("M132").Value = -10295.6519
$ws.Range("N132").Value = -30576.8
$ws.Range("H141").Value = 7215.7915
$ws.Range("I141").Value = 1698.1765
$ws.Range("J141").Value = 20615.715
$ws.Range("K141").Value = 5094.529500000001
$ws.Range("L141").Value = 61847.145
$ws.Range("M141").Value = 85.47049999999945
$ws.Range("N141").Value = -72207.145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1333.3334
$ws.Range("I45").Value = 1400
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1400
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -1023
$ws.Range("N45").Value = -1754
$ws.Range("H109").Value = 15673.8
$ws.Range("J109").Value = 15673.8
$ws.Range("L109").Value = 15673.8
$ws.Range("N109").Value = -18447.8
$ws.Range("H132").Value = 1586741.4
$ws.Range("I132").Value = 4397.317
$ws.Range("J132").Value = 3294007.5
$ws.Range("K132").Value = 13191.951
$ws.Range("L132").Value = 9882022.5
$ws.Range("M132").Value = -10661.951
$ws.Range("N132").Value = -9887082.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1891.4546
$ws.Range("I107").Value = 1886.2858
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1886.2858
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 33.71419999999989
$ws.Range("N107").Value = -5840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 531.6923
$ws.Range("I22").Value = 249.33333
$ws.Range("J22").Value = 773.7143
$ws.Range("K22").Value = 249.33333
$ws.Range("L22").Value = 773.7143
$ws.Range("M22").Value = 100.66667
$ws.Range("N22").Value = -1473.7143
$ws.Range("H33").Value = 154
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H62").Value = 27475
$ws.Range("I62").Value = 3300
$ws.Range("K62").Value = 3300
$ws.Range("M62").Value = -2676
$ws.Range("H65").Value = 27475
$ws.Range("I65").Value = 3300
$ws.Range("K65").Value = 16500
$ws.Range("M65").Value = -13380
$ws.Range("H122").Value = 71429896
$ws.Range("I122").Value = 90909850
$ws.Range("J122").Value = 3404.6667
$ws.Range("K122").Value = 272729550
$ws.Range("L122").Value = 10214.0001
$ws.Range("M122").Value = -272727100
$ws.Range("N122").Value = -15114.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1019.5
$ws.Range("I17").Value = 481.66666
$ws.Range("J17").Value = 1250
$ws.Range("K17").Value = 1444.99998
$ws.Range("L17").Value = 3750
$ws.Range("M17").Value = -1275.99998
$ws.Range("N17").Value = -4088
$ws.Range("H46").Value = 733.5714
$ws.Range("I46").Value = 689.1667
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 2067.5001
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -1976.5001
$ws.Range("N46").Value = -3182

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 898.7
$ws.Range("I9").Value = 568
$ws.Range("K9").Value = 568
$ws.Range("M9").Value = -398
$ws.Range("H95").Value = 22428.572
$ws.Range("J95").Value = 22428.572
$ws.Range("L95").Value = 22428.572
$ws.Range("N95").Value = -27920.572
$ws.Range("H126").Value = 5959300.5
$ws.Range("I126").Value = 11370015
$ws.Range("K126").Value = 34110045
$ws.Range("M126").Value = -34107575
$ws.Range("H132").Value = 3712.9565
$ws.Range("I132").Value = 2182.4546
$ws.Range("J132").Value = 5115.9165
$ws.Range("K132").Value = 6547.3638
$ws.Range("L132").Value = 15347.7495
$ws.Range("M132").Value = -4017.3638
$ws.Range("N132").Value = -20407.7495

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2675.8333
$ws.Range("I7").Value = 2675.8333
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2675.8333
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2563.8333
$ws.Range("N7").ClearContents()
$ws.Range("H100").Value = 2411.7666
$ws.Range("I100").Value = 1706.25
$ws.Range("J100").Value = 3218.0715
$ws.Range("K100").Value = 1706.25
$ws.Range("L100").Value = 3218.0715
$ws.Range("M100").Value = -1165.25
$ws.Range("N100").Value = -4300.0715
$ws.Range("H111").Value = 29462.334
$ws.Range("J111").Value = 29462.334
$ws.Range("L111").Value = 29462.334
$ws.Range("N111").Value = -37642.334
$ws.Range("H126").Value = 2675.8333
$ws.Range("I126").Value = 2675.8333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8027.499899999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5557.499899999999
$ws.Range("N126").ClearContents()
